# Actualización automática 2025-07-31 08:55:10
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Narrow column E slightly (engine adds a fixed 5/6 offset when
# converting ColumnWidth to the stored character width, so compensate)
$ws.Columns.Item(5).ColumnWidth = 22.166666666666668

# Row 2 (OTROS)
$ws.Range("D2").Value = 6786.72
$ws.Range("E2").Value = -6786.72

# Row 3 (PORCELANATO)
$ws.Range("D3").Value = 14328.13
$ws.Range("E3").Value = -604.7899999999991
$ws.Range("F3").Value = 1.044070175336325

# Row 4 (TOTAL)
$ws.Range("D4").Value = 21114.85
$ws.Range("E4").Value = -7391.509999999999
$ws.Range("F4").Value = 1.53860867689644
